$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Approved bookings")

$newRows = @(
    @{ A = "gh";            B = "eroha@gmail.com";      C = "2022-10-10 00:00:00 UTC"; D = "2022-10-13 00:00:00 UTC" },
    @{ A = "name";          B = "email@email.com";      C = "2022-10-16 00:00:00 UTC"; D = "2022-10-17 00:00:00 UTC" },
    @{ A = "Ivan Knyazev";  B = "pjesuss120@gmail.com"; C = "2022-12-10 00:00:00 UTC"; D = "2022-12-15 00:00:00 UTC" }
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
}
